# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Swap the displayed country names for two row pairs (Suiza/Japon, Malta/Guinea Ecuatorial)
# - Refresh the per-country COVID figures (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose data changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp on row 1 -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 12:59"

# --- Country-name swaps ------------------------------------------------------
# Row 51 used to read "Japon" / row 52 used to read "Suiza" -> they swap.
$ws.Range("A51").Value = "Suiza"
$ws.Range("A52").Value = "Japon"

# Row 136 used to read "Guinea Ecuatorial" / row 137 used to read "Malta" -> they swap.
$ws.Range("A136").Value = "Malta"
$ws.Range("A137").Value = "Guinea Ecuatorial"

# --- Updated figures ----------------------------------------------------------
# Row 16 (Iran)
$ws.Range("B16").Value = 550757
$ws.Range("C16").Value = 5471
$ws.Range("D16").Value = 442674
$ws.Range("E16").Value = 76433
$ws.Range("G16").Value = 304
$ws.Range("H16").Value = 31650

# Row 34 (Rumania)
$ws.Range("B34").Value = 196004
$ws.Range("C34").Value = 4902
$ws.Range("D34").Value = 141089
$ws.Range("E34").Value = 48752
$ws.Range("G34").Value = 98
$ws.Range("H34").Value = 6163

# Row 37 (Nepal)
$ws.Range("B37").Value = 148509
$ws.Range("C37").Value = 3637
$ws.Range("D37").Value = 102820
$ws.Range("E37").Value = 44877
$ws.Range("G37").Value = 21
$ws.Range("H37").Value = 812

# Row 42 (Emiratos Arabes Unidos)
$ws.Range("B42").Value = 120710
$ws.Range("C42").Value = 1578
$ws.Range("D42").Value = 113364
$ws.Range("E42").Value = 6872
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 474

# Row 51 (now Suiza)
$ws.Range("B51").Value = 97019
$ws.Range("C51").Value = 5256
$ws.Range("D51").Value = 55700
$ws.Range("E51").Value = 39275
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 2044

# Row 52 (now Japon)
$ws.Range("B52").Value = 93933
$ws.Range("D52").Value = 87107
$ws.Range("E52").Value = 5147
$ws.Range("H52").Value = 1679

# Row 91 (Malasia)
$ws.Range("B91").Value = 23804
$ws.Range("C91").Value = 847
$ws.Range("D91").Value = 15417
$ws.Range("E91").Value = 8183
$ws.Range("G91").Value = 5
$ws.Range("H91").Value = 204

# Row 101 (Senegal)
$ws.Range("B101").Value = 15508
$ws.Range("C101").Value = 24
$ws.Range("D101").Value = 14026
$ws.Range("E101").Value = 1161

# Row 131 (Hong Kong)
$ws.Range("B131").Value = 5281
$ws.Range("C131").Value = 11
$ws.Range("D131").Value = 5019
$ws.Range("E131").Value = 157

# Row 136 (now Malta)
$ws.Range("B136").Value = 5137
$ws.Range("C136").Value = 111
$ws.Range("D136").Value = 3384
$ws.Range("E136").Value = 1704
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = 49

# Row 137 (now Guinea Ecuatorial)
$ws.Range("B137").Value = 5074
$ws.Range("D137").Value = 4954
$ws.Range("E137").Value = 37
$ws.Range("H137").Value = 83

# Row 175 (Gibraltar)
$ws.Range("B175").Value = 630
$ws.Range("C175").Value = 9
$ws.Range("D175").Value = 495
$ws.Range("E175").Value = 135
